# Apply edits described by the target diff for "Chemin critique.xlsx"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New comments on E11/F11 (Visualiser une requete), E13/F13 (Visualiser les dossiers), E14/F14 (Visualiser un dossier) ---
$txt_E11 = @'
Visualiser une requête: (En dev.)
La visualiton d'une requête à eu un bon avancement cette semaine. L'affichage en loi même est terminé comme nous le désirons et il ne restera que le relié avec l'API et corrigé quelques petits problèmes qui peut surgir.
Le temps prévu pour cette tâche cette semaine était d'environ 3h, ainsi nous sommes restés dans les alentours du temps prévues. Nous pensions que le visuel allait être un peu plus complexe à réaliser, mais pour la visualisation d'une requête, le visuel a été relativement simple à faire.
'@
$ws.Range("E11").AddComment($txt_E11) | Out-Null

$txt_F11 = @'
Visualiser une requête: (Terminé)
Avec l'arrivé de l'API qui était fonctionnel ainsi que les résultats d'analyse que nous devions faire, la visualiser d'une requête est devenu beaucoup plus complexe que nous l'avions initialement prévue. Même si nous avons réussi à le terminé, le cas d'utilisation à prit beaucoup plus de temps que nous l'avons prévus à cause de léger problème avec l'API, mais surtout à cause du nombre de type d'analyse que nous devions mettre dans l'application.
Le temps prévus à la tâche était d'environ 4h30 cette semaine, ce qui a été largement dépassé de 2h15. Un grand excès qui peut être vu, comme dit plus haut, par la sous-estimation du nombre d'analyse qui était requis pour les requêtes.
'@
$ws.Range("F11").AddComment($txt_F11) | Out-Null

$txt_E13 = @'
Visualiser les dossiers: (En dev.)
Étant donné que ce cas d'utilisation visuel était le premier qui devait être fait, nous avions prévus que celui-ci allait prendre beaucoup de temps cette semaine et le temps utilisé de 11h30 à respecter cette échéancier tout en permettant de l'avancée à un bon rythme.
Nous avions prévus travaillé sur ce cas d'utilisation environ 11h à 12h cette semaine, un temps qui a été respecté tout en finissant le visuel de la liste.
'@
$ws.Range("E13").AddComment($txt_E13) | Out-Null

$txt_F13 = @'
Visualiser les dossiers: (Terminé)
Le cas d'utilisation a été terminé dans un temps très raisonnable tout en étant ce que nous nous attendions à avoir.
Nous avions prévus travaillé environ 3h sur ce cas d'utilisation pour faire l'appel à l'API et les problèmes qui pouvait subvenir avec cette appel. Même si nous avons réussis à terminé le cas, sa finition à été clairement plus petite que ce que nous avons prévus (1h45 de moins). La sous-estimation de temps peut être vu à cause du nombre de temps que l'API nous as prit dans les autres tâches, ce qui nous a influencé à estimé une plus grande durée pour ce cas d'utilisation.
'@
$ws.Range("F13").AddComment($txt_F13) | Out-Null

$txt_E14 = @'
Visualiser un dossier: (En dev.)
Le cas d'utilisation, même si son avancement est bon, nous a prit trop de temps pour ce que nous nous attendions à prendre.
Nous avons estimé que le temps requis pour ce cas cette semaine serait de 4h, un temps qui a été dépassé de 30 minutes. Cette excès de temps peut être venu du fait que nous apprenions le language utilisé et que cet apprentisage à prit plus de temps que nous l'avions initialement prévus.
'@
$ws.Range("E14").AddComment($txt_E14) | Out-Null

$txt_F14 = @'
Visualiser un dossier: (Terminé)
Ce cas d'utilisation à prit beaucoup plus de temps que nous l'avons initialement prévu et même s'il est terminé, nous n'avons passé trop de temps à le finir.
Le temps prévus pour effectué cette tâche cette semaine était d'environ 1h30, un temps qui a été plus que doublé. Cet excès de temps peut être remarqué à cause de la difficulté de l'API qui a influencé plusieurs tâches de l'application.
'@
$ws.Range("F14").AddComment($txt_F14) | Out-Null

# --- 2. Wording fix in the "Creer le dossier: (Termine)" comment (F12) ---
$f12_text = @'
Créer le dossier: (Terminé)
L'avancement fait cette semaine sur la création de dossier était beaucoup plus productif cette semaine. Le formulaire a pu être terminé et même si le temps a été dépassé 2h45, l'avancement fait a été plus rapide qu'auparavant.

De cette manière, l'excès de temps peut être vu à cause des problèmes avec l'accès au API qui sont arrivés dans ce cas d'utilisation. Le bouton pour la création reste la cause principale de l'excès de temps
'@
$ws.Range("F12").Comment.Text($f12_text)

# --- 3. Fill in the placeholder text of the "API Test: (Termine)" comment (G19) ---
$g19_text = @'
API Test: (Terminé)
Les tests de l'API ont été terminés facilement et efficacement cette semaine. 
Le temps estimé de 1h pour terminé le test restant des requêtes à été respecté ce qui a permis à ce test d'être dans les temps que nous avons prévus initialement.
'@
$ws.Range("G19").Comment.Text($g19_text)

# --- 4. Fill in previously-empty G25/H25 and G27/H27 with 0, copying the
#        "value entered" cell formatting (fill + border) from row 16, which
#        already has the same G/H style pattern used once a value is entered. ---
$ws.Range("G16").Copy()
$ws.Range("G25").PasteSpecial(-4122)
$ws.Range("H16").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("G16").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("H16").Copy()
$ws.Range("H27").PasteSpecial(-4122)

$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0

# --- 5. Update the sheet view: drop the frozen/scrolled topLeftCell and move
#        the active selection from J21 to K17. ---
$ws.Range("K17").Select()

